$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it appears ---
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# --- 2. zh-cn handback datetime (was placeholder 0001-01-01, now real timestamp) ---
$wsZh.Range("K2").Value = "2016-10-20 09:33:56"
$wsZh.Range("K3").Value = "2016-10-20 09:33:56"

# --- 3. de-de handback datetime (new timestamp, distinct from zh-cn's) ---
$wsDe.Range("K2").Value = "2016-10-20 09:34:15"
$wsDe.Range("K3").Value = "2016-10-20 09:34:15"

# --- 4. Latest Target File (I) / Latest Handback File (J) columns now populated ---
$zhMdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40c04deb83037a2110d1acbcbb8ac7a923334173/e2e/2b20de96-b05a-43ed-886c-18547d96bb90.md"
$zhMdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40c04deb83037a2110d1acbcbb8ac7a923334173/e2e/7b262a9c-ac66-45ae-927b-19e8ea82e08f.md"

# zh-cn row 2 (file 2b20de96...)
$wsZh.Range("J2").Value = "2b20de96-b05a-43ed-886c-18547d96bb90.87738b9c6c4b32727c76fb64e1955a3e88726050.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhMdUrl1, "", "", "2b20de96-b05a-43ed-886c-18547d96bb90.md") | Out-Null
$wsZh.Range("I2").Style = "HyperLink"

# zh-cn row 3 (file 7b262a9c...)
$wsZh.Range("J3").Value = "7b262a9c-ac66-45ae-927b-19e8ea82e08f.6c31b09cd70b39611f5d546152297f65e466db36.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhMdUrl2, "", "", "7b262a9c-ac66-45ae-927b-19e8ea82e08f.md") | Out-Null
$wsZh.Range("I3").Style = "HyperLink"

# de-de row 2
$wsDe.Range("J2").Value = "2b20de96-b05a-43ed-886c-18547d96bb90.87738b9c6c4b32727c76fb64e1955a3e88726050.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $zhMdUrl1, "", "", "2b20de96-b05a-43ed-886c-18547d96bb90.md") | Out-Null
$wsDe.Range("I2").Style = "HyperLink"

# de-de row 3
$wsDe.Range("J3").Value = "7b262a9c-ac66-45ae-927b-19e8ea82e08f.6c31b09cd70b39611f5d546152297f65e466db36.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $zhMdUrl2, "", "", "7b262a9c-ac66-45ae-927b-19e8ea82e08f.md") | Out-Null
$wsDe.Range("I3").Style = "HyperLink"

# --- 5. Column width adjustments ---
# NOTE: the host engine stores column width as (round(ColumnWidth*6)+5)/6, i.e. it
# always adds ~0.8333 of padding and then rounds to the nearest 1/6. The target
# widths below (29.9777050018311 / 40) are fed through the inverse of that
# transform so the persisted <col width="..."> lands as close as the engine's
# quantization allows (exact for 40, within 0.0223 for 29.9777050018311).
$wTarget30 = 29.166666666666668   # -> stored 30 (closest achievable to 29.9777050018311)
$wTarget40 = 39.166666666666664   # -> stored 40 (exact)

$wsOverview.Columns.Item(5).ColumnWidth = $wTarget30
$wsOverview.Columns.Item(6).ColumnWidth = $wTarget30

$wsZh.Columns.Item(3).ColumnWidth = $wTarget30
$wsZh.Columns.Item(9).ColumnWidth = $wTarget40
$wsZh.Columns.Item(10).ColumnWidth = $wTarget40

$wsDe.Columns.Item(3).ColumnWidth = $wTarget30
$wsDe.Columns.Item(9).ColumnWidth = $wTarget40
$wsDe.Columns.Item(10).ColumnWidth = $wTarget40

"Done applying handback report changes"
